$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.538380445877948
$ws.Cells.Item(2, 3).Value = 0.04233200773009571
$ws.Cells.Item(2, 4).Value = 0.6652258404945997
$ws.Cells.Item(2, 5).Value = 0.2711893627522741
$ws.Cells.Item(2, 7).Value = 0.002518792601187466
$ws.Cells.Item(2, 10).Value = 0.1410305611426637
$ws.Cells.Item(2, 11).Value = 0.5032047446278796
$ws.Cells.Item(2, 15).Value = 5.892821856924797
$ws.Cells.Item(3, 2).Value = 0.5028528678126349
$ws.Cells.Item(3, 3).Value = 0.03748746243596202
$ws.Cells.Item(3, 4).Value = 0.6552871200087793
$ws.Cells.Item(3, 5).Value = 0.2662916400292659
$ws.Cells.Item(3, 7).Value = 0.002522075926011473
$ws.Cells.Item(3, 10).Value = 0.1377689756088287
$ws.Cells.Item(3, 11).Value = 0.465401726843794
$ws.Cells.Item(3, 15).Value = 5.899319320728324
$ws.Cells.Item(4, 2).Value = 0.4812742419334768
$ws.Cells.Item(4, 3).Value = 0.0345163390492047
$ws.Cells.Item(4, 4).Value = 0.6495256094539172
$ws.Cells.Item(4, 5).Value = 0.2634298749776178
$ws.Cells.Item(4, 7).Value = 0.002524199042694553
$ws.Cells.Item(4, 10).Value = 0.1358465485516831
$ws.Cells.Item(4, 11).Value = 0.4423871120722822
$ws.Cells.Item(4, 15).Value = 5.906085425792554
$ws.Cells.Item(5, 2).Value = 0.4725403297707089
$ws.Cells.Item(5, 3).Value = 0.03330649273696906
$ws.Cells.Item(5, 4).Value = 0.6472635146630807
$ws.Cells.Item(5, 5).Value = 0.2623002565817103
$ws.Cells.Item(5, 7).Value = 0.002525091255716247
$ws.Cells.Item(5, 10).Value = 0.135083291272224
$ws.Cells.Item(5, 11).Value = 0.4330582644384435
$ws.Cells.Item(5, 15).Value = 5.90954040744603
$ws.Cells.Item(6, 2).Value = 0.4710936774963557
$ws.Cells.Item(6, 3).Value = 0.03310565498350115
$ws.Cells.Item(6, 4).Value = 0.646893077053079
$ws.Cells.Item(6, 5).Value = 0.2621148932554647
$ws.Cells.Item(6, 7).Value = 0.002525241041770289
$ws.Cells.Item(6, 10).Value = 0.1349577692757151
$ws.Cells.Item(6, 11).Value = 0.4315122346056626
$ws.Cells.Item(6, 15).Value = 5.910156232638059
$ws.Cells.Item(7, 2).Value = 0.481156211748754
$ws.Cells.Item(7, 3).Value = 0.03450001889298449
$ws.Cells.Item(7, 4).Value = 0.649494754722042
$ws.Cells.Item(7, 5).Value = 0.2634144924633119
$ws.Cells.Item(7, 7).Value = 0.002524210965790865
$ws.Cells.Item(7, 10).Value = 0.1358361734462008
$ws.Cells.Item(7, 11).Value = 0.4422610978720911
$ws.Cells.Item(7, 15).Value = 5.906129196408557
$ws.Cells.Item(8, 2).Value = 0.5260818640995524
$ws.Cells.Item(8, 3).Value = 0.04066091763100133
$ws.Cells.Item(8, 4).Value = 0.6617282365491519
$ws.Cells.Item(8, 5).Value = 0.2694704347873582
$ws.Cells.Item(8, 7).Value = 0.002519902504900995
$ws.Cells.Item(8, 10).Value = 0.1398893078218393
$ws.Cells.Item(8, 11).Value = 0.4901296638275312
$ws.Cells.Item(8, 15).Value = 5.894485596734796
$ws.Cells.Item(9, 2).Value = 0.616039014601597
$ws.Cells.Item(9, 3).Value = 0.05276849066387967
$ws.Cells.Item(9, 4).Value = 0.6884231641164433
$ws.Cells.Item(9, 5).Value = 0.2825013471280755
$ws.Cells.Item(9, 7).Value = 0.002512299914125139
$ws.Cells.Item(9, 10).Value = 0.1484753268403054
$ws.Cells.Item(9, 11).Value = 0.5855491617182338
$ws.Cells.Item(9, 15).Value = 5.893710926350366
$ws.Cells.Item(10, 2).Value = 0.6832572697290971
$ws.Cells.Item(10, 3).Value = 0.06167904637788979
$ws.Cells.Item(10, 4).Value = 0.7096884576464788
$ws.Cells.Item(10, 5).Value = 0.292782377403249
$ws.Cells.Item(10, 7).Value = 0.002507224832058141
$ws.Cells.Item(10, 10).Value = 0.1551752504531549
$ws.Cells.Item(10, 11).Value = 0.6565922492789582
$ws.Cells.Item(10, 15).Value = 5.906635107262616
$ws.Cells.Item(11, 2).Value = 0.7140805198862097
$ws.Cells.Item(11, 3).Value = 0.06573588057106861
$ws.Cells.Item(11, 4).Value = 0.7197223785069014
$ws.Cells.Item(11, 5).Value = 0.2976137859902295
$ws.Cells.Item(11, 7).Value = 0.002505025759848191
$ws.Cells.Item(11, 10).Value = 0.1583089795384609
$ws.Cells.Item(11, 11).Value = 0.6891146103224628
$ws.Cells.Item(11, 15).Value = 5.91545514288606
$ws.Cells.Item(12, 2).Value = 0.7257875376953962
$ws.Cells.Item(12, 3).Value = 0.06727255640403484
$ws.Cells.Item(12, 4).Value = 0.7235737793768067
$ws.Cells.Item(12, 5).Value = 0.2994655646333015
$ws.Cells.Item(12, 7).Value = 0.002504208703032723
$ws.Cells.Item(12, 10).Value = 0.1595080315109669
$ws.Cells.Item(12, 11).Value = 0.7014591600774054
$ws.Cells.Item(12, 15).Value = 5.919218655401664
$ws.Cells.Item(13, 2).Value = 0.7232646744080284
$ws.Cells.Item(13, 3).Value = 0.06694158687226093
$ws.Cells.Item(13, 4).Value = 0.7227420091149384
$ws.Cells.Item(13, 5).Value = 0.2990657621645596
$ws.Cells.Item(13, 7).Value = 0.002504383974490869
$ws.Cells.Item(13, 10).Value = 0.1592492433620407
$ws.Cells.Item(13, 11).Value = 0.6987992548065449
$ws.Cells.Item(13, 15).Value = 5.918389266999327
$ws.Cells.Item(14, 2).Value = 0.7150429654240043
$ws.Cells.Item(14, 3).Value = 0.06586229523720988
$ws.Cells.Item(14, 4).Value = 0.7200381981268436
$ws.Cells.Item(14, 5).Value = 0.2977656873720065
$ws.Cells.Item(14, 7).Value = 0.002504958226257599
$ws.Cells.Item(14, 10).Value = 0.1584073780755517
$ws.Cells.Item(14, 11).Value = 0.6901296240085344
$ws.Cells.Item(14, 15).Value = 5.915756276723016
$ws.Cells.Item(15, 2).Value = 0.7100114686431027
$ws.Cells.Item(15, 3).Value = 0.06520125402850852
$ws.Cells.Item(15, 4).Value = 0.718388777360019
$ws.Cells.Item(15, 5).Value = 0.2969722492038187
$ws.Cells.Item(15, 7).Value = 0.002505312012301777
$ws.Cells.Item(15, 10).Value = 0.1578933238322833
$ws.Cells.Item(15, 11).Value = 0.6848229954484282
$ws.Cells.Item(15, 15).Value = 5.914198673913347
$ws.Cells.Item(16, 2).Value = 0.681247803954335
$ws.Cells.Item(16, 3).Value = 0.06141398683080013
$ws.Cells.Item(16, 4).Value = 0.7090399636777818
$ws.Cells.Item(16, 5).Value = 0.2924697435792538
$ws.Cells.Item(16, 7).Value = 0.0025073707457065
$ws.Cells.Item(16, 10).Value = 0.1549721845012044
$ws.Cells.Item(16, 11).Value = 0.6544709199085617
$ws.Cells.Item(16, 15).Value = 5.906117930731654
$ws.Cells.Item(17, 2).Value = 0.6636648072723972
$ws.Cells.Item(17, 3).Value = 0.0590914539441485
$ws.Cells.Item(17, 4).Value = 0.7033970129693046
$ws.Cells.Item(17, 5).Value = 0.2897471843608415
$ws.Cells.Item(17, 7).Value = 0.002508661729984732
$ws.Cells.Item(17, 10).Value = 0.1532021785279909
$ws.Cells.Item(17, 11).Value = 0.6359030183190839
$ws.Cells.Item(17, 15).Value = 5.90191433750266
$ws.Cells.Item(18, 2).Value = 0.6535746451653495
$ws.Cells.Item(18, 3).Value = 0.05775591501864596
$ws.Cells.Item(18, 4).Value = 0.700185245338389
$ws.Cells.Item(18, 5).Value = 0.2881957835010382
$ws.Cells.Item(18, 7).Value = 0.002509414591506202
$ws.Cells.Item(18, 10).Value = 0.1521922017664821
$ws.Cells.Item(18, 11).Value = 0.6252425231463121
$ws.Cells.Item(18, 15).Value = 5.899773299286807
$ws.Cells.Item(19, 2).Value = 0.6501622701637473
$ws.Cells.Item(19, 3).Value = 0.05730378092434307
$ws.Cells.Item(19, 4).Value = 0.699103620675686
$ws.Cells.Item(19, 5).Value = 0.2876730033436132
$ws.Cells.Item(19, 7).Value = 0.002509671272649006
$ws.Cells.Item(19, 10).Value = 0.1518516286594007
$ws.Cells.Item(19, 11).Value = 0.6216363837699816
$ws.Cells.Item(19, 15).Value = 5.899095893311937
$ws.Cells.Item(20, 2).Value = 0.665534157818513
$ws.Cells.Item(20, 3).Value = 0.05933865868652788
$ws.Cells.Item(20, 4).Value = 0.7039942055261577
$ws.Cells.Item(20, 5).Value = 0.2900355002031674
$ws.Cells.Item(20, 7).Value = 0.002508523234892038
$ws.Cells.Item(20, 10).Value = 0.1533897618613906
$ws.Cells.Item(20, 11).Value = 0.6378776113904223
$ws.Cells.Item(20, 15).Value = 5.902333168896519
$ws.Cells.Item(21, 2).Value = 0.7174569356965037
$ws.Cells.Item(21, 3).Value = 0.06617929774431275
$ws.Cells.Item(21, 4).Value = 0.7208309680748926
$ws.Cells.Item(21, 5).Value = 0.2981469474815768
$ws.Cells.Item(21, 7).Value = 0.00250478912960558
$ws.Cells.Item(21, 10).Value = 0.1586543182188791
$ws.Cells.Item(21, 11).Value = 0.6926753195857032
$ws.Cells.Item(21, 15).Value = 5.916518149160936
$ws.Cells.Item(22, 2).Value = 0.7515948619698065
$ws.Cells.Item(22, 3).Value = 0.07065259969098747
$ws.Cells.Item(22, 4).Value = 0.7321365250778058
$ws.Cells.Item(22, 5).Value = 0.3035778201402195
$ws.Cells.Item(22, 7).Value = 0.002502440060161172
$ws.Cells.Item(22, 10).Value = 0.1621671630391717
$ws.Cells.Item(22, 11).Value = 0.7286578738384435
$ws.Cells.Item(22, 15).Value = 5.928257982892944
$ws.Cells.Item(23, 2).Value = 0.7333563146842721
$ws.Cells.Item(23, 3).Value = 0.06826489566397242
$ws.Cells.Item(23, 4).Value = 0.7260749344955286
$ws.Cells.Item(23, 5).Value = 0.3006674018742714
$ws.Cells.Item(23, 7).Value = 0.002503685466650076
$ws.Cells.Item(23, 10).Value = 0.1602856814931783
$ws.Cells.Item(23, 11).Value = 0.7094379572360765
$ws.Cells.Item(23, 15).Value = 5.921766061905316
$ws.Cells.Item(24, 2).Value = 0.664688966670326
$ws.Cells.Item(24, 3).Value = 0.05922689830377692
$ws.Cells.Item(24, 4).Value = 0.7037241137248884
$ws.Cells.Item(24, 5).Value = 0.2899051095196654
$ws.Cells.Item(24, 7).Value = 0.002508585815337151
$ws.Cells.Item(24, 10).Value = 0.1533049316939241
$ws.Cells.Item(24, 11).Value = 0.6369848528617297
$ws.Cells.Item(24, 15).Value = 5.902142956561278
$ws.Cells.Item(25, 2).Value = 0.5915049376159516
$ws.Cells.Item(25, 3).Value = 0.04949036330576462
$ws.Cells.Item(25, 4).Value = 0.6809115294444155
$ws.Cells.Item(25, 5).Value = 0.2788521609619821
$ws.Cells.Item(25, 7).Value = 0.002514266573351096
$ws.Cells.Item(25, 10).Value = 0.146083994465684
$ws.Cells.Item(25, 11).Value = 0.5595704427652493
$ws.Cells.Item(25, 15).Value = 5.891554259591885
